# alerts_cdb.xlsx — "Adding new kind of graphs & adding a new line for the excel"
#
# The underlying CDB alert percentages were regenerated (new raw data pulled
# in upstream), so every existing data point B2:J12 moves to a new value,
# two previously-blank cells (G5, H7) now carry a value, the B2:J12 number
# format gains two decimals, and a new "Bilan" (balance) summary row (13) is
# appended below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Recalculated values for the existing data block (rows 2-12) -------
$ws.Range("B2").Value = -31.682522530926871
$ws.Range("D2").Value = 32.962910532276332
$ws.Range("F2").Value = -48.761222212613539
$ws.Range("G2").Value = 99.521055975372292
$ws.Range("H2").Value = -139.68958475153161
$ws.Range("I2").Value = -69.73951598004804
$ws.Range("J2").Value = -49.440361998596508

$ws.Range("B3").Value = -18.325301969017811
$ws.Range("C3").Value = -28.7133331986164
$ws.Range("D3").Value = 23.591383137643621

$ws.Range("B4").Value = -78.202867617092338
$ws.Range("C4").Value = -52.485894376605437
$ws.Range("D4").Value = -48.006444348933549
$ws.Range("F4").Value = 76.01605949048664
$ws.Range("G4").Value = 52.987598647125147
$ws.Range("H4").Value = 69.699348640539739
$ws.Range("J4").Value = -9.8227828683494014

$ws.Range("B5").Value = -46.596465277662588
$ws.Range("D5").Value = -14.76669363559925
$ws.Range("F5").Value = -77.48459756252069
$ws.Range("G5").Value = -99.865995160936365
$ws.Range("H5").Value = -75.599056276928422
$ws.Range("J5").Value = 67.279314369178138

$ws.Range("B6").Value = -21.613501696916959
$ws.Range("C6").Value = -111.4687301192356
$ws.Range("D6").Value = -16.879889436914389
$ws.Range("E6").Value = 75.954689048688763
$ws.Range("F6").Value = -60.352870592787369
$ws.Range("G6").Value = 5.8583673382575974
$ws.Range("H6").Value = -35.497112208929657
$ws.Range("I6").Value = 27.08324630359084

$ws.Range("B7").Value = -63.495560582948038
$ws.Range("C7").Value = -104.2319487500237
$ws.Range("D7").Value = -51.440507320548193
$ws.Range("E7").Value = -108.4023799815524
$ws.Range("F7").Value = -103.04618667897159
$ws.Range("G7").Value = -110.1318309629569
$ws.Range("H7").Value = -100
$ws.Range("J7").Value = -66.074062166323529

$ws.Range("B8").Value = 81.144740769845072
$ws.Range("D8").Value = -48.220590025565073
$ws.Range("E8").Value = -49.605559493861477
$ws.Range("F8").Value = -50.509977827051003
$ws.Range("G8").Value = 108.7812812751876
$ws.Range("H8").Value = -60.708923490440228
$ws.Range("J8").Value = -10.32495130526906

$ws.Range("B9").Value = -26.44493040807787
$ws.Range("C9").Value = -54.635576654565767
$ws.Range("D9").Value = -55.648505178364807
$ws.Range("F9").Value = -1.418217511369148
$ws.Range("G9").Value = -61.993743382999469
$ws.Range("H9").Value = -70.33815179557314
$ws.Range("J9").Value = -79.18344267642901

$ws.Range("B10").Value = 43.095801532699163
$ws.Range("D10").Value = -33.080500916282382
$ws.Range("G10").Value = -29.051191284092351
$ws.Range("J10").Value = 33.803853616183922

$ws.Range("B11").Value = -12.05584468197825
$ws.Range("D11").Value = -30.128961194669991
$ws.Range("F11").Value = -9.5144525527172696
$ws.Range("G11").Value = -20.263210546206349
$ws.Range("H11").Value = -26.46824092765598
$ws.Range("I11").Value = 10.876388267778241

$ws.Range("B12").Value = 19.419042495965581
$ws.Range("D12").Value = 12.784234461849421
$ws.Range("E12").Value = -35.163435452921327
$ws.Range("F12").Value = 6.1283185840707954
$ws.Range("G12").Value = -0.31746031746031739
$ws.Range("H12").Value = 26.41960549910344
$ws.Range("I12").Value = -6.1902831663674744

# --- 2. Number format for the data block goes from "0" to "0.00" ----------
$ws.Range("B2:J12").NumberFormat = "0.00"

# --- 3. New "Bilan" (balance) row 13, matching row 12's look & feel -------
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B12:J12").Copy()
$ws.Range("B13:J13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A13").Value = "Bilan"
$ws.Range("B13").Value = -169.38207490470049
$ws.Range("C13").Value = -3351.5969524299312
$ws.Range("D13").Value = -1146.4077416878761
$ws.Range("E13").Value = -3036.7466324437169
$ws.Range("F13").Value = -2520.9576606668788
$ws.Range("G13").Value = 527.16561371143644
$ws.Range("H13").Value = -3299.1082223560052
$ws.Range("I13").Value = -3790.841674027015
$ws.Range("J13").Value = -1347.826128814796

# --- 4. Selection moves to N9 (matches the saved sheet view) --------------
[void]$ws.Range("N9").Select()
